# Update the "合肥-漫展信息" workbook per the upstream data refresh.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: style a freshly-inserted index cell (column A) so it matches the
# bold / centered / thin-bordered look used by every other row in these
# sheets (style index 1 in the original file).
# ---------------------------------------------------------------------------
function Set-IndexCellStyle($cell) {
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1
}

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibitions) - sheet index 1
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F2").Value = 1684
$ws1.Range("F3").Value = 9129
$ws1.Range("F4").Value = 115
$ws1.Range("F5").Value = 508
$ws1.Range("F6").Value = 708
$ws1.Range("F7").Value = 1377
$ws1.Range("F8").Value = 204
$ws1.Range("F9").Value = 61
$ws1.Range("F10").Value = 98
$ws1.Range("F11").Value = 5922
$ws1.Range("F13").Value = 388
$ws1.Range("F15").Value = 4506
$ws1.Range("F16").Value = 15
$ws1.Range("F18").Value = 1152
$ws1.Range("F21").Value = 30

# Insert the new event row (长丰·莓可可游戏动漫展) above the old row 22,
# pushing everything from row 22 down by one.
$ws1.Rows.Item(22).Insert()
Set-IndexCellStyle $ws1.Cells.Item(22,1)

$ws1.Cells.Item(22,1).Value = 21
$ws1.Cells.Item(22,2).Value = "2024-07-27"
$ws1.Cells.Item(22,3).Value = "长丰·莓可可游戏动漫展"
$ws1.Cells.Item(22,4).Value = "长寿路12号 长丰宾馆·梅山饭店(长寿路店)"
$ws1.Cells.Item(22,5).Value = "2024.07.27 10:00-07.27 17:00"
$ws1.Cells.Item(22,6).Value = 0
$ws1.Cells.Item(22,7).Value = 40
$ws1.Cells.Item(22,8).Value = "https://show.bilibili.com/platform/detail.html?id=87796"
$ws1.Cells.Item(22,9).Value = "//i2.hdslb.com/bfs/openplatform/202406/MLTfeikq1718823574810.png"

# The rows that shifted down one position also had their "want to go" counts
# (column F) refreshed.
$ws1.Range("F23").Value = 261
$ws1.Range("F25").Value = 2771

# ---------------------------------------------------------------------------
# Sheet "演出" (Performances) - sheet index 2
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 38

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All types) - sheet index 4
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F2").Value = 1684
$ws4.Range("F3").Value = 9129
$ws4.Range("F4").Value = 115
$ws4.Range("F5").Value = 38
$ws4.Range("F6").Value = 508
$ws4.Range("F7").Value = 708
$ws4.Range("F8").Value = 1377
$ws4.Range("F9").Value = 204
$ws4.Range("F10").Value = 61
$ws4.Range("F11").Value = 98
$ws4.Range("F12").Value = 5922
$ws4.Range("F14").Value = 388
$ws4.Range("F16").Value = 4506
$ws4.Range("F17").Value = 15
$ws4.Range("F19").Value = 1152
$ws4.Range("F22").Value = 30

# Insert the new event row (长丰·莓可可游戏动漫展) above the old row 23,
# pushing everything from row 23 down by one.
$ws4.Rows.Item(23).Insert()
Set-IndexCellStyle $ws4.Cells.Item(23,1)

$ws4.Cells.Item(23,1).Value = 22
$ws4.Cells.Item(23,2).Value = "2024-07-27"
$ws4.Cells.Item(23,3).Value = "长丰·莓可可游戏动漫展"
$ws4.Cells.Item(23,4).Value = "长寿路12号 长丰宾馆·梅山饭店(长寿路店)"
$ws4.Cells.Item(23,5).Value = "2024.07.27 10:00-07.27 17:00"
$ws4.Cells.Item(23,6).Value = 0
$ws4.Cells.Item(23,7).Value = 40
$ws4.Cells.Item(23,8).Value = "https://show.bilibili.com/platform/detail.html?id=87796"
$ws4.Cells.Item(23,9).Value = "//i2.hdslb.com/bfs/openplatform/202406/MLTfeikq1718823574810.png"

$ws4.Range("F24").Value = 261
$ws4.Range("F26").Value = 2771
